# Add a "State" column to the demographic table and update values.
#
# The table originally has 2 columns (Subgroup, Percentage of Students).
# We need 3 columns (Subgroup, Davidson County, State), with updated /
# reordered values, and every data cell left-aligned ("jc val=left")
# instead of the original right alignment.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$W_NS = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# "Compact" is a paragraph style referenced throughout the document but it
# is not actually defined in styles.xml. Range.InsertXML validates pStyle
# references against the known style table and silently drops unknown
# ones, so we register a temporary placeholder style, use it while we
# rebuild the cells (so pStyle references survive InsertXML), and then
# remove the placeholder again before saving so styles.xml stays exactly
# as it was.
$tempStyle = $d.Styles.Add("Compact", 1)

# Final desired content, by row: label (col1), Davidson County value (col2),
# State value (col3).
$finalRows = @(
  @("Subgroup", "Davidson County", "State"),
  @("Black/Hispanic/Native American", "65.2%", "32.9%"),
  @("Economically Disadvantaged", "75.3%", "57.9%"),
  @("English Learners", "16.2%", "4.6%"),
  @("Students with Disabilities", "12.4%", "14%")
)

# Insert a brand new blank column. In this object model, Columns.Add()
# with no arguments always inserts the new column at the front: the
# existing "Subgroup" column shifts from 1->2, and the existing
# "Percentage of Students" column shifts from 2->3.
$t.Columns.Add() | Out-Null

$rowCount = $t.Rows.Count
for ($r = 1; $r -le $rowCount; $r++) {
  $label = $finalRows[$r - 1][0]
  $davidson = $finalRows[$r - 1][1]
  $state = $finalRows[$r - 1][2]

  $col1 = $t.Cell($r, 1)
  $col2 = $t.Cell($r, 2)
  $col3 = $t.Cell($r, 3)

  $xml1 = "<w:p $W_NS><w:pPr><w:pStyle w:val=`"Compact`"/><w:jc w:val=`"left`"/></w:pPr><w:r><w:t xml:space=`"preserve`">$label</w:t></w:r></w:p>"
  $xml2 = "<w:p $W_NS><w:pPr><w:pStyle w:val=`"Compact`"/><w:jc w:val=`"left`"/></w:pPr><w:r><w:t xml:space=`"preserve`">$davidson</w:t></w:r></w:p>"
  $xml3 = "<w:p $W_NS><w:pPr><w:pStyle w:val=`"Compact`"/><w:jc w:val=`"left`"/></w:pPr><w:r><w:t xml:space=`"preserve`">$state</w:t></w:r></w:p>"

  $col1.Range.InsertXML($xml1)
  $col2.Range.InsertXML($xml2)
  $col3.Range.InsertXML($xml3)

  if ($r -eq 1) {
    # Header row: restore the bottom border + bottom vertical alignment
    # on column 1 (lost because Columns.Add() creates a plain blank cell)
    # and set up the same formatting on the new column 3 cell.
    $col1.Borders.Item(-3).LineStyle = 1
    $col1.VerticalAlignment = 3

    $col3.Borders.Item(-3).LineStyle = 1
    $col3.VerticalAlignment = 3
  }
}

# Clean up the temporary style registration so styles.xml is left
# untouched by our edit.
$d.Styles.Item("Compact").Delete()
